$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column F ("Reviewed By") before the old "Description of the Change"
# column, which shifts the old F column to G.
$ws.Columns("F:F").Insert()
$ws.Columns("F:F").ColumnWidth = $ws.Columns("E:E").ColumnWidth

# New log entries (rows 4 and 5), populated in the same order the original
# author entered them so that shared-string indices line up.
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 43145
$ws.Range("C4").Value = 0.78125
$ws.Range("D4").Value = "Team_04_M1_D3_Communication_Policy"
$ws.Range("G4").Value = "Initial Version"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 43146
$ws.Range("C5").Value = 0.78194444444444444
$ws.Range("D5").Value = "Team_04_M1_D5_Liaison"
$ws.Range("G5").Value = "Initial Version"

$ws.Range("E4").Value = "2100-Prerana"
$ws.Range("E5").Value = "2100-Prerana"

# New header for the inserted column
$ws.Range("F1").Value = "Reviewed By"

# Update selection to reflect where the user left off editing
$ws.Range("F5").Select()
